# Regenerate merged AHB files
# - Rename the "_old"/"_new" suffixed header labels to "_FV2310"/"_FV2404"
# - Turn the header range A1:U66 into a proper Excel Table ("Table1")
# - Freeze the header row (pane split after row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the column headers in row 1 -------------------------------
$headers = @(
  "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
  "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310",
  "diff",
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
  "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Convert the used range A1:U66 into an Excel Table -----------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row (split below row 1) -------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Edit applied: headers renamed, table created, top row frozen."
